# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets
# to reflect newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 84
    $ws.Range("F4").Value = 1482
    $ws.Range("F5").Value = 17
    $ws.Range("F6").Value = 32
    $ws.Range("F7").Value = 117
    $ws.Range("F9").Value = 268
}
